# Add a new "Quantity On Purchase" column (E) to the stock-enquiry report,
# populated with the purchase-order quantity for every product row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in E1, matching the style used by the other header cells
# (bold, centered, thin border) by copying D1's formatting.
$ws.Range("E1").Value = "Quantity On Purchase"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Quantity-on-purchase values for every data row (row -> value).
$purchaseQty = @{
    2  = 0;     3  = 0;     4  = 0;     5  = 0;     6  = 0;
    7  = 0;     8  = 6500;  9  = 36500; 10 = 0;     11 = 0;
    12 = 0;     13 = 0;     14 = 0;     15 = 0;     16 = 0;
    17 = 6500;  18 = 36500; 19 = 3285;  20 = 0;
    22 = 0;     23 = 0;     24 = 0;     25 = 0;     26 = 0;
    27 = 0;     28 = 6500;  29 = 36500; 30 = 2000;  31 = 0;
    32 = 0;     33 = 0;     34 = 0;     35 = 0;     36 = 0;
    37 = 6500;  38 = 36500; 39 = 2000;  40 = 0;
    42 = 0;     43 = 0;     44 = 0;     45 = 0;     46 = 0;
    48 = 0;     49 = 0;     50 = 0;     51 = 0;     52 = 87000;
    53 = 0;     54 = 0;
    56 = 0;     57 = 2150;  58 = 0;     59 = 0;      60 = 6000;
    61 = 0;     62 = 5000;  63 = 0;     64 = 0
}

foreach ($row in $purchaseQty.Keys) {
    $ws.Cells.Item($row, 5).Value = $purchaseQty[$row]
}

# Re-fit the column widths now that column E (and its longer header) exists,
# mirroring the "best fit" widths Excel computes for this data.
$ws.Columns.Item(1).ColumnWidth = 16.17
$ws.Columns.Item(2).ColumnWidth = 48.17
$ws.Columns.Item(3).ColumnWidth = 16.02
$ws.Columns.Item(4).ColumnWidth = 16.02
$ws.Columns.Item(5).ColumnWidth = 19.59

# Leave the cursor where the author's last save left it.
$ws.Range("H10").Select() | Out-Null

